# "example of a fight completed"
# - Questy: quest #3 (Kill hobgoblin) EndStage 3 -> 1
# - Dialogi: mark stages 8 & 20 as completed (B column 0 -> 1)
# - Dialogi: collapse two multi-line "Me:" choice texts into single-line
#   versions (with the row heights that go with the shorter text)
# - Dialogi: add a new "Delete" flag column (L) with True/False markers
# - Dialogi: add a full example fight dialogue block (QuestID 333, rows 21-24)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Dialogi" - new example fight dialogue (QuestID 333)
# Cell text is entered in this particular order so that newly introduced
# shared strings land at the expected table positions.
# ---------------------------------------------------------------------------
$dialogi = $wb.Worksheets.Item("Dialogi")

$dialogi.Range("A23").Value = 333
$dialogi.Range("B23").Value = 3
$dialogi.Range("C23").Value = 0
$dialogi.Range("D23").Value = "False"
$dialogi.Range("E23").Value = 3
$dialogi.Range("H23").Value = "Blastoise:    I'll catch You next time !"
$dialogi.Range("K23").Value = 2

$dialogi.Range("A21").Value = 333
$dialogi.Range("B21").Value = 3
$dialogi.Range("C21").Value = 0
$dialogi.Range("D21").Value = "True"
$dialogi.Range("E21").Value = 1
$dialogi.Range("F21").Value = 3
$dialogi.Range("H21").Value = "Me:    2. Run away !    3. Defend yourself !"
$dialogi.Range("H21").WrapText = $true
$dialogi.Range("J21").Value = "True"
$dialogi.Range("K21").Value = 1

# Collapse the multi-line choice text into a single wrapped line, and let
# the rows reflow to their natural (shorter) height
$dialogi.Range("H12").Value = "Me:    1. Cool, I like cash    2. [CHA] That quest was hard, reward should be bigger    3. I'll just kill you and take more cash"
$dialogi.Rows.Item(12).RowHeight = 30

$dialogi.Range("H5").Value = "Me:    1. I will handle it!    2. I do not have time for this right now."
$dialogi.Rows.Item(5).AutoFit() | Out-Null

$dialogi.Range("A24").Value = 333
$dialogi.Range("B24").Value = 3
$dialogi.Range("C24").Value = 0
$dialogi.Range("D24").Value = "False"
$dialogi.Range("E24").Value = 4
$dialogi.Range("H24").Value = "Blastoise:    I'll kill Ya !"
$dialogi.Range("K24").Value = 3

# New "Delete" marker column
$dialogi.Range("L1").Value = "Delete"
$dialogi.Columns.Item(12).ColumnWidth = 18.14

$dialogi.Range("L7").Value = "False"
$dialogi.Range("L14").Value = "True"
$dialogi.Range("L20").Value = "False"

# Blank (but text-formatted / wrapped) placeholder cell for the middle of
# the fight example
$dialogi.Range("A22").Value = 333
$dialogi.Range("B22").Value = 3
$dialogi.Range("C22").Value = 0
$dialogi.Range("D22").Value = "False"
$dialogi.Range("E22").Value = 2
$dialogi.Range("H22").Value = "'"
$dialogi.Range("H22").WrapText = $true
$dialogi.Range("H22").Value = ""
$dialogi.Range("K22").Value = 1

# Flip the two "completed" flags
$dialogi.Range("B8").Value = 1
$dialogi.Range("B20").Value = 1

# ---------------------------------------------------------------------------
# Sheet "Questy"
# ---------------------------------------------------------------------------
$questy = $wb.Worksheets.Item("Questy")

$questy.Range("C5").Value = 1

$questy.Activate()
$questy.Range("D10").Select() | Out-Null

$dialogi.Activate()
$dialogi.Range("H22").Select() | Out-Null
